$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (best achievable values given the runtime's internal
# rounding of ColumnWidth to the nearest 1/6 character unit)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.833333333333332   # B -> ~12.73046875
$ws.Columns.Item(3).ColumnWidth = 12.0                 # C -> ~12.86328125
$ws.Columns.Item(4).ColumnWidth = 25.333333333333336   # D -> ~26.1328125
$ws.Columns.Item(5).ColumnWidth = 12.666666666666668   # E -> ~13.46484375
$ws.Columns.Item(6).ColumnWidth = 21.166666666666664   # F -> ~21.9296875
$ws.Columns.Item(7).ColumnWidth = 25.5                 # G -> ~26.3984375

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14.65
$ws.Rows.Item(3).RowHeight = 14.65
$ws.Rows.Item(4).RowHeight = 14.65
$ws.Rows.Item(7).RowHeight = 55.9
$ws.Rows.Item(8).RowHeight = 48.85
$ws.Rows.Item(9).RowHeight = 55.25
$ws.Rows.Item(10).RowHeight = 57.4
$ws.Rows.Item(13).RowHeight = 52.9

# ---------------------------------------------------------------------------
# Fill in the newly completed "Rectangle" test cases (rows 7-13)
# ---------------------------------------------------------------------------

# Row 7 - __init__ / successful creation
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'color =  "red",                                                         length =  4,                                                                      width = 5'
$ws.Range("G7").Value = "The rectangle instance is created successfully with the attributes correctly set."

# Row 8 - __init__ / blank color
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'color =  "  ",                                                         length =  4,                                                                      width = 5'
$ws.Range("G8").Value = "Value error"

# Row 9 - __init__ / length not an integer
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'color =  "red",                                                         length = "length",                                                                     width = 5'
$ws.Range("G9").Value = "Value error"

# Row 10 - __init__ / width not an integer
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'color =  "red",                                                         length =  4,                                                                      width =  "width"'
$ws.Range("G10").Value = "Value error"

# Row 11 - __str__
$ws.Range("E11").Value = 'color =  "red",                                                         length =  4,                                                                      width = 5'
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "`"The shape color is red.\nThis rectangle has four sides with the lengths of 4, 5, 4 and 5 centimeters.`""

# Row 12 - calculate_area
$ws.Range("E12").Value = 'color =  "red",                                                         length =  4,                                                                      width = 5'
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = 20

# Row 13 - calculate_perimeter (E13/F13 pick up the same "filled in" look
# as the other Method Inputs / Preconditions cells, copied from E7's format)
$ws.Range("E7").Copy()
$ws.Range("E13:F13").PasteSpecial(-4122)
$ws.Range("E13").Value = 'color =  "red",                                                         length =  4,                                                                      width = 5'
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = 18

# ---------------------------------------------------------------------------
# Selection / scroll position
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
